$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Remove the four "project idea" slides that are no longer relevant,
#    keeping only the "Hallucination Detection" and "Roadmap" slides.
# ---------------------------------------------------------------------------
$null = $p.Slides.Item(1).Delete()
$null = $p.Slides.Item(1).Delete()
$null = $p.Slides.Item(1).Delete()
$null = $p.Slides.Item(1).Delete()

# ---------------------------------------------------------------------------
# 2. First remaining slide: drop the leading "5. " numbering from the title.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape1 = $slide1.Shapes.Item(1)
$titleShape1.TextFrame.TextRange.Text = "Hallucination Detection & Mitigation for LLMs"

# ---------------------------------------------------------------------------
# 3. Second remaining slide: rename "General Roadmap" -> "Planned Roadmap"
#    and append an extra note after the Month 6 bullet.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$titleShape2 = $slide2.Shapes.Item(1)
$titleShape2.TextFrame.TextRange.Text = "Planned Roadmap"

$bodyShape2 = $slide2.Shapes.Item(2)
$tr2 = $bodyShape2.TextFrame.TextRange
$null = $tr2.InsertAfter([char]13 + "An additional 1 month will be reserved for any unexpected event.")
$paras2 = $tr2.Paragraphs()
$lastPara2 = $tr2.Paragraphs($paras2.Count, 1)
$lastPara2.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# 4. Bump the cached "today" date shown in the footer date placeholder
#    (slide master + every slide layout) from 12/9/2025 to 15/9/2025.
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "12/9/2025") {
                $sh.TextFrame.TextRange.Text = "15/9/2025"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes
for ($k = 1; $k -le $master.CustomLayouts.Count; $k++) {
    $layout = $master.CustomLayouts.Item($k)
    Update-DateShape $layout.Shapes
}

Write-Output "Edit complete. Slides remaining: $($p.Slides.Count)"
